$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-refresh) values for every cell the crypto-price refresh touched.
# Cell refs map 1:1 onto the diff hunks in the commit (price/volume updates plus
# the handful of rank swaps where two rows traded places).
$updates = @(
    @{ Cell = 'D2'; Text = '63.852.47' }
    @{ Cell = 'E2'; Text = '  -4.98%  ' }
    @{ Cell = 'D3'; Text = '3.301.98' }
    @{ Cell = 'E3'; Text = '  -5.61%  ' }
    @{ Cell = 'E4'; Text = '  +0.27%  ' }
    @{ Cell = 'D5'; Text = '179.09' }
    @{ Cell = 'E5'; Text = '  -10.74%  ' }
    @{ Cell = 'D6'; Text = '526.63' }
    @{ Cell = 'E6'; Text = '  -4.33%  ' }
    @{ Cell = 'D7'; Text = '0.603' }
    @{ Cell = 'E7'; Text = '  -0.45%  ' }
    @{ Cell = 'D8'; Text = '3.297.57' }
    @{ Cell = 'E8'; Text = '  -5.50%  ' }
    @{ Cell = 'E9'; Text = '  +0.24%  ' }
    @{ Cell = 'D10'; Text = '0.610' }
    @{ Cell = 'E10'; Text = '  -6.44%  ' }
    @{ Cell = 'D11'; Text = '58.00' }
    @{ Cell = 'E11'; Text = '  -8.04%  ' }
    @{ Cell = 'D12'; Text = '0.133' }
    @{ Cell = 'E12'; Text = '  -6.95%  ' }
    @{ Cell = 'D13'; Text = '0.0000259' }
    @{ Cell = 'E13'; Text = '  -4.07%  ' }
    @{ Cell = 'D14'; Text = '9.11' }
    @{ Cell = 'E14'; Text = '  -6.98%  ' }
    @{ Cell = 'D15'; Text = '3.845.73' }
    @{ Cell = 'E15'; Text = '  -4.93%  ' }
    @{ Cell = 'D16'; Text = '3.312.26' }
    @{ Cell = 'E16'; Text = '  -4.89%  ' }
    @{ Cell = 'D17'; Text = '0.117' }
    @{ Cell = 'E17'; Text = '  -5.14%  ' }
    @{ Cell = 'D18'; Text = '63.997.81' }
    @{ Cell = 'E18'; Text = '  -4.30%  ' }
    @{ Cell = 'D19'; Text = '17.46' }
    @{ Cell = 'E19'; Text = '  -4.59%  ' }
    @{ Cell = 'D20'; Text = '11.13' }
    @{ Cell = 'E20'; Text = '  -5.49%  ' }
    @{ Cell = 'D21'; Text = '0.958' }
    @{ Cell = 'E21'; Text = '  -6.35%  ' }
    @{ Cell = 'D22'; Text = '373.96' }
    @{ Cell = 'E22'; Text = '  -4.29%  ' }
    @{ Cell = 'D23'; Text = '3.79' }
    @{ Cell = 'E23'; Text = '  -4.97%  ' }
    @{ Cell = 'B24'; Text = 'Litecoin' }
    @{ Cell = 'C24'; Text = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' }
    @{ Cell = 'D24'; Text = '80.74' }
    @{ Cell = 'E24'; Text = '  -2.10%  ' }
    @{ Cell = 'B25'; Text = 'RenderToken' }
    @{ Cell = 'C25'; Text = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D25'; Text = '11.13' }
    @{ Cell = 'E25'; Text = '  -10.91%  ' }
    @{ Cell = 'D26'; Text = '3.92' }
    @{ Cell = 'E26'; Text = '  +1.17%  ' }
    @{ Cell = 'D27'; Text = '6.08' }
    @{ Cell = 'E27'; Text = '  -1.41%  ' }
    @{ Cell = 'D28'; Text = '2.68' }
    @{ Cell = 'E28'; Text = '  -4.66%  ' }
    @{ Cell = 'D29'; Text = '11.42' }
    @{ Cell = 'E29'; Text = '  -6.30%  ' }
    @{ Cell = 'D30'; Text = '8.37' }
    @{ Cell = 'E30'; Text = '  -5.00%  ' }
    @{ Cell = 'D31'; Text = '28.92' }
    @{ Cell = 'E31'; Text = '  -6.44%  ' }
    @{ Cell = 'D32'; Text = '646.56' }
    @{ Cell = 'E32'; Text = '  -4.69%  ' }
    @{ Cell = 'D33'; Text = '6.64' }
    @{ Cell = 'E33'; Text = '  -4.30%  ' }
    @{ Cell = 'D34'; Text = '11.27' }
    @{ Cell = 'E34'; Text = '  -3.51%  ' }
    @{ Cell = 'D35'; Text = '0.106' }
    @{ Cell = 'E35'; Text = '  -4.75%  ' }
    @{ Cell = 'D36'; Text = '59.32' }
    @{ Cell = 'E36'; Text = '  -6.82%  ' }
    @{ Cell = 'E37'; Text = '  -0.08%  ' }
    @{ Cell = 'D38'; Text = '0.390' }
    @{ Cell = 'E38'; Text = '  -1.29%  ' }
    @{ Cell = 'D39'; Text = '36.68' }
    @{ Cell = 'E39'; Text = '  -5.21%  ' }
    @{ Cell = 'E40'; Text = '  +0.46%  ' }
    @{ Cell = 'D41'; Text = '0.0₃0701' }
    @{ Cell = 'E41'; Text = '  +4.64%  ' }
    @{ Cell = 'B42'; Text = 'Maker' }
    @{ Cell = 'C42'; Text = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' }
    @{ Cell = 'D42'; Text = '2.949.63' }
    @{ Cell = 'E42'; Text = '  -3.57%  ' }
    @{ Cell = 'B43'; Text = 'Kaspa' }
    @{ Cell = 'C43'; Text = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' }
    @{ Cell = 'D43'; Text = '0.125' }
    @{ Cell = 'E43'; Text = '  -4.54%  ' }
    @{ Cell = 'D44'; Text = '2.48' }
    @{ Cell = 'E44'; Text = '  -4.15%  ' }
    @{ Cell = 'D45'; Text = '2.70' }
    @{ Cell = 'E45'; Text = '  -8.77%  ' }
    @{ Cell = 'D46'; Text = '0.0398' }
    @{ Cell = 'E46'; Text = '  +0.03%  ' }
    @{ Cell = 'D47'; Text = '2.65' }
    @{ Cell = 'E47'; Text = '  -4.43%  ' }
    @{ Cell = 'B48'; Text = 'ApeXProtocol' }
    @{ Cell = 'C48'; Text = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex' }
    @{ Cell = 'D48'; Text = '3.09' }
    @{ Cell = 'E48'; Text = '  +8.13%  ' }
    @{ Cell = 'B49'; Text = 'Stacks' }
    @{ Cell = 'C49'; Text = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' }
    @{ Cell = 'D49'; Text = '2.79' }
    @{ Cell = 'E49'; Text = '  +6.55%  ' }
    @{ Cell = 'D50'; Text = '0.126' }
    @{ Cell = 'E50'; Text = '  -0.44%  ' }
    @{ Cell = 'B51'; Text = 'Monero' }
    @{ Cell = 'C51'; Text = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D51'; Text = '135.27' }
    @{ Cell = 'E51'; Text = '  -1.35%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Cells such as "179.09" or "58.00" parse as pure numbers, and a plain
    # .Value assignment would let Excel coerce them to numeric cells. Force
    # text storage, then drop back to the Normal style so no stray number
    # format lingers on the cell (matches the unstyled source cells).
    if ($u.Text -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Text
    }
}
